$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.582.23"
$ws.Range("E2").Value = "  +0.64%  "
# Row 3
$ws.Range("D3").Value = "2.637.86"
$ws.Range("E3").Value = "  +0.72%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.82%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
# Row 7
$ws.Range("E7").Value = "  +0.00%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.547"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.84%  "
# Row 9
$ws.Range("D9").Value = "2.637.45"
$ws.Range("E9").Value = "  +0.74%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.19%  "
# Row 11
$ws.Range("E11").Value = "  +0.70%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.07%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.78%  "
# Row 15
$ws.Range("D15").Value = "3.112.77"
$ws.Range("E15").Value = "  +0.61%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000184"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.34%  "
# Row 17
$ws.Range("D17").Value = "67.566.83"
$ws.Range("E17").Value = "  +0.77%  "
# Row 18
$ws.Range("D18").Value = "2.613.39"
$ws.Range("E18").Value = "  -0.11%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.36%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "366.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.51%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.66%  "
# Row 22
$ws.Range("E22").Value = "  -0.46%  "
# Row 23
$ws.Range("E23").Value = "  +7.03%  "
# Row 24
$ws.Range("E24").Value = "  -0.01%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.08%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.67%  "
# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000105"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.48%  "
# Row 28
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.756.03"
$ws.Range("E28").Value = "  +0.01%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "583.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.41%  "
# Row 30
$ws.Range("E30").Value = "  -0.40%  "
# Row 31
$ws.Range("E31").Value = "  -2.74%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "
# Row 33
$ws.Range("E33").Value = "  -0.14%  "
# Row 34
$ws.Range("E34").Value = "  -2.22%  "
# Row 35
$ws.Range("E35").Value = "  +0.06%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.52%  "
# Row 38
$ws.Range("E38").Value = "  +2.10%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.372"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.51%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.07%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.69%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "
# Row 48
$ws.Range("E48").Value = "  -4.40%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.07%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.62%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.628"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.37%  "
